$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look numeric but must stay as literal text
# (to preserve exact formatting like trailing zeros / grouping dots),
# so force text format before assignment.
$textForceCells = @("D4", "D5", "D6", "D10", "D11", "D12", "D14", "D15", "D19", "D20", "D22", "D23", "D24", "D26", "D28", "D30", "D31", "D32", "D33", "D40", "D42", "D44", "D45", "D46", "D48", "D49")
foreach ($addr in $textForceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated values

# Row 2
$ws.Range("D2").Value = '42.904.97'
$ws.Range("E2").Value = '  +1.44%  '

# Row 3
$ws.Range("D3").Value = '2.288.26'
$ws.Range("E3").Value = '  -0.52%  '

# Row 4
$ws.Range("D4").Value = '1.01'
$ws.Range("E4").Value = '  +0.40%  '

# Row 5
$ws.Range("D5").Value = '314.35'
$ws.Range("E5").Value = '  -0.66%  '

# Row 6
$ws.Range("D6").Value = '105.26'
$ws.Range("E6").Value = '  +2.51%  '

# Row 7
$ws.Range("E7").Value = '  +0.60%  '

# Row 8
$ws.Range("E8").Value = '  +0.16%  '

# Row 9
$ws.Range("E9").Value = '  +0.25%  '

# Row 10
$ws.Range("D10").Value = '39.71'
$ws.Range("E10").Value = '  +0.51%  '

# Row 11
$ws.Range("D11").Value = '0.0906'
$ws.Range("E11").Value = '  -0.23%  '

# Row 12
$ws.Range("D12").Value = '8.40'
$ws.Range("E12").Value = '  +0.17%  '

# Row 13
$ws.Range("E13").Value = '  +2.59%  '

# Row 14
$ws.Range("D14").Value = '0.995'
$ws.Range("E14").Value = '  +3.51%  '

# Row 15
$ws.Range("D15").Value = '15.26'
$ws.Range("E15").Value = '  +0.41%  '

# Row 16
$ws.Range("D16").Value = '2.634.73'
$ws.Range("E16").Value = '  -0.56%  '

# Row 17
$ws.Range("D17").Value = '2.331.79'
$ws.Range("E17").Value = '  +1.77%  '

# Row 18
$ws.Range("D18").Value = '42.795.93'
$ws.Range("E18").Value = '  +0.91%  '

# Row 19
$ws.Range("D19").Value = '7.44'
$ws.Range("E19").Value = '  -0.12%  '

# Row 20
$ws.Range("D20").Value = '13.82'
$ws.Range("E20").Value = '  +21.73%  '

# Row 21
$ws.Range("E21").Value = '  -0.09%  '

# Row 22
$ws.Range("D22").Value = '74.01'

# Row 23
$ws.Range("D23").Value = '3.59'
$ws.Range("E23").Value = '  +1.48%  '

# Row 24
$ws.Range("D24").Value = '265.74'
$ws.Range("E24").Value = '  -3.85%  '

# Row 25
$ws.Range("E25").Value = '  -1.85%  '

# Row 26
$ws.Range("D26").Value = '1.00'
$ws.Range("E26").Value = '  +0.24%  '

# Row 27
$ws.Range("E27").Value = '  +0.06%  '

# Row 28
$ws.Range("D28").Value = '7.22'
$ws.Range("E28").Value = '  +23.12%  '

# Row 29
$ws.Range("E29").Value = '  -0.27%  '

# Row 30
$ws.Range("D30").Value = '22.53'
$ws.Range("E30").Value = '  -0.92%  '

# Row 31
$ws.Range("D31").Value = '37.27'
$ws.Range("E31").Value = '  -0.66%  '

# Row 32
$ws.Range("D32").Value = '167.31'
$ws.Range("E32").Value = '  +0.96%  '

# Row 33
$ws.Range("D33").Value = '0.0878'
$ws.Range("E33").Value = '  +0.44%  '

# Row 34
$ws.Range("E34").Value = '  -3.08%  '

# Row 35
$ws.Range("E35").Value = '  -0.89%  '

# Row 36
$ws.Range("E36").Value = '  -4.09%  '

# Row 37
$ws.Range("E37").Value = '  -0.16%  '

# Row 38
$ws.Range("E38").Value = '  -3.78%  '

# Row 39
$ws.Range("E39").Value = '  +2.98%  '

# Row 40
$ws.Range("D40").Value = '2.67'
$ws.Range("E40").Value = '  -3.58%  '

# Row 41
$ws.Range("E41").Value = '  +4.72%  '

# Row 42
$ws.Range("D42").Value = '70.95'
$ws.Range("E42").Value = '  +1.83%  '

# Row 43
$ws.Range("E43").Value = '  +3.03%  '

# Row 44
$ws.Range("B44").Value = 'BitcoinSV'
$ws.Range("C44").Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
$ws.Range("D44").Value = '94.71'
$ws.Range("E44").Value = '  -0.01%  '

# Row 45
$ws.Range("B45").Value = 'FirstDigitalUSD'
$ws.Range("C45").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D45").Value = '1.01'
$ws.Range("E45").Value = '  +0.15%  '

# Row 46
$ws.Range("D46").Value = '12.16'
$ws.Range("E46").Value = '  +1.13%  '

# Row 47
$ws.Range("D47").Value = '1.740.49'
$ws.Range("E47").Value = '  +9.55%  '

# Row 48
$ws.Range("D48").Value = '113.01'
$ws.Range("E48").Value = '  -0.03%  '

# Row 49
$ws.Range("D49").Value = '80.08'
$ws.Range("E49").Value = '  -1.01%  '

# Row 50
$ws.Range("E50").Value = '  -0.67%  '

# Row 51
$ws.Range("E51").Value = '  -2.41%  '
